# Daily auto-update of "variazione giornaliera" data: the data window rolls forward by one day.
# Each existing row 2-10 takes on the values that were previously in the row below it (rows 3-11),
# and a brand-new row 11 is appended for the newest date. Two data points (Lazio and P.A. Trento)
# were not yet available for the newest date, so those two cells are left blank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 44400
$ws.Range("B2").Value = 0.27
$ws.Range("C2").Value = 0.25
$ws.Range("D2").Value = 0.26
$ws.Range("F2").Value = 0.51
$ws.Range("G2").Value = 0.22
$ws.Range("H2").Value = 0.75
$ws.Range("I2").Value = 0.45
$ws.Range("J2").Value = 0.32
$ws.Range("K2").Value = 0.33
$ws.Range("L2").Value = 0.2
$ws.Range("M2").Value = 0.32
$ws.Range("N2").Value = 0.48
$ws.Range("O2").Value = 0.18
$ws.Range("P2").Value = 0.18
$ws.Range("Q2").Value = 0.97
$ws.Range("R2").Value = 0.75
$ws.Range("S2").Value = 0.66
$ws.Range("T2").Value = 0.5600000000000001
$ws.Range("U2").Value = 0.22
$ws.Range("V2").Value = 0.78

# Row 3
$ws.Range("A3").Value = 44401
$ws.Range("B3").Value = 0.28
$ws.Range("C3").Value = 0.28
$ws.Range("D3").Value = 0.3
$ws.Range("E3").Value = 0.33
$ws.Range("F3").Value = 0.5600000000000001
$ws.Range("G3").Value = 0.25
$ws.Range("H3").Value = 0.8
$ws.Range("I3").Value = 0.51
$ws.Range("J3").Value = 0.33
$ws.Range("K3").Value = 0.36
$ws.Range("L3").Value = 0.18
$ws.Range("M3").Value = 0.28
$ws.Range("N3").Value = 0.5
$ws.Range("O3").Value = 0.2
$ws.Range("P3").Value = 0.19
$ws.Range("Q3").Value = 1.01
$ws.Range("R3").Value = 0.79
$ws.Range("S3").Value = 0.72
$ws.Range("T3").Value = 0.63
$ws.Range("U3").Value = 0.25
$ws.Range("V3").Value = 0.8

# Row 4
$ws.Range("A4").Value = 44402
$ws.Range("B4").Value = 0.31
$ws.Range("C4").Value = 0.27
$ws.Range("D4").Value = 0.33
$ws.Range("E4").Value = 0.34
$ws.Range("F4").Value = 0.63
$ws.Range("G4").Value = 0.27
$ws.Range("H4").Value = 0.83
$ws.Range("J4").Value = 0.34
$ws.Range("K4").Value = 0.37
$ws.Range("L4").Value = 0.19
$ws.Range("M4").Value = 0.3
$ws.Range("N4").Value = 0.51
$ws.Range("P4").Value = 0.22
$ws.Range("Q4").Value = 1.07
$ws.Range("R4").Value = 0.82
$ws.Range("S4").Value = 0.79
$ws.Range("T4").Value = 0.66
$ws.Range("U4").Value = 0.27
$ws.Range("V4").Value = 0.82

# Row 5
$ws.Range("A5").Value = 44403
$ws.Range("B5").Value = 0.3
$ws.Range("D5").Value = 0.35
$ws.Range("E5").Value = 0.33
$ws.Range("F5").Value = 0.6899999999999999
$ws.Range("G5").Value = 0.28
$ws.Range("H5").Value = 0.84
$ws.Range("I5").Value = 0.53
$ws.Range("J5").Value = 0.33
$ws.Range("K5").Value = 0.39
$ws.Range("L5").Value = 0.18
$ws.Range("M5").Value = 0.29
$ws.Range("N5").Value = 0.53
$ws.Range("O5").Value = 0.21
$ws.Range("P5").Value = 0.21
$ws.Range("Q5").Value = 1.14
$ws.Range("R5").Value = 0.85
$ws.Range("S5").Value = 0.84
$ws.Range("T5").Value = 0.65
$ws.Range("V5").Value = 0.85

# Row 6
$ws.Range("A6").Value = 44404
$ws.Range("B6").Value = 0.28
$ws.Range("D6").Value = 0.39
$ws.Range("E6").Value = 0.32
$ws.Range("F6").Value = 0.7
$ws.Range("G6").Value = 0.35
$ws.Range("H6").Value = 0.8
$ws.Range("I6").Value = 0.5600000000000001
$ws.Range("J6").Value = 0.34
$ws.Range("M6").Value = 0.28
$ws.Range("N6").Value = 0.54
$ws.Range("O6").Value = 0.23
$ws.Range("Q6").Value = 1.2
$ws.Range("R6").Value = 0.82
$ws.Range("S6").Value = 0.86
$ws.Range("T6").Value = 0.6899999999999999
$ws.Range("U6").Value = 0.29
$ws.Range("V6").Value = 0.88

# Row 7
$ws.Range("A7").Value = 44405
$ws.Range("B7").Value = 0.34
$ws.Range("C7").Value = 0.24
$ws.Range("D7").Value = 0.43
$ws.Range("E7").Value = 0.33
$ws.Range("F7").Value = 0.72
$ws.Range("G7").Value = 0.33
$ws.Range("H7").Value = 0.8100000000000001
$ws.Range("I7").Value = 0.58
$ws.Range("J7").Value = 0.35
$ws.Range("K7").Value = 0.44
$ws.Range("N7").Value = 0.5
$ws.Range("P7").Value = 0.23
$ws.Range("Q7").Value = 1.27
$ws.Range("R7").Value = 0.83
$ws.Range("S7").Value = 0.9399999999999999
$ws.Range("T7").Value = 0.76
$ws.Range("U7").Value = 0.3
$ws.Range("V7").Value = 0.89

# Row 8
$ws.Range("A8").Value = 44406
$ws.Range("C8").Value = 0.26
$ws.Range("D8").Value = 0.49
$ws.Range("G8").Value = 0.36
$ws.Range("H8").Value = 0.8
$ws.Range("I8").Value = 0.5600000000000001
$ws.Range("J8").Value = 0.36
$ws.Range("M8").Value = 0.33
$ws.Range("O8").Value = 0.26
$ws.Range("P8").Value = 0.22
$ws.Range("Q8").Value = 1.36
$ws.Range("R8").Value = 0.86
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 0.8100000000000001
$ws.Range("V8").Value = 0.87

# Row 9
$ws.Range("A9").Value = 44407
$ws.Range("B9").Value = 0.37
$ws.Range("C9").Value = 0.3
$ws.Range("D9").Value = 0.48
$ws.Range("E9").Value = 0.34
$ws.Range("F9").Value = 0.73
$ws.Range("G9").Value = 0.39
$ws.Range("H9").Value = 0.78
$ws.Range("I9").Value = 0.55
$ws.Range("K9").Value = 0.46
$ws.Range("L9").Value = 0.17
$ws.Range("M9").Value = 0.32
$ws.Range("N9").Value = 0.46
$ws.Range("O9").Value = 0.29
$ws.Range("P9").Value = 0.25
$ws.Range("Q9").Value = 1.38
$ws.Range("R9").Value = 0.91
$ws.Range("S9").Value = 1.07
$ws.Range("T9").Value = 0.84
$ws.Range("U9").Value = 0.25
$ws.Range("V9").Value = 0.93

# Row 10
$ws.Range("A10").Value = 44408
$ws.Range("B10").Value = 0.39
$ws.Range("C10").Value = 0.36
$ws.Range("D10").Value = 0.47
$ws.Range("E10").Value = 0.36
$ws.Range("F10").Value = 0.77
$ws.Range("G10").Value = 0.4
$ws.Range("H10").Value = 0.74
$ws.Range("I10").Value = 0.5600000000000001
$ws.Range("J10").Value = 0.38
$ws.Range("K10").Value = 0.5
$ws.Range("M10").Value = 0.34
$ws.Range("N10").Value = 0.45
$ws.Range("O10").Value = 0.3
$ws.Range("P10").Value = 0.27
$ws.Range("Q10").Value = 1.39
$ws.Range("R10").Value = 0.96
$ws.Range("S10").Value = 1.11
$ws.Range("T10").Value = 0.86
$ws.Range("U10").Value = 0.23
$ws.Range("V10").Value = 0.9399999999999999

# Row 11
$ws.Range("A11").Value = 44409
$ws.Range("B11").Value = 0.37
$ws.Range("C11").Value = 0.39
$ws.Range("D11").Value = 0.49
$ws.Range("E11").Value = 0.37
$ws.Range("G11").Value = 0.41
$ws.Range("H11").ClearContents()
$ws.Range("I11").Value = 0.61
$ws.Range("K11").Value = 0.52
$ws.Range("L11").Value = 0.19
$ws.Range("M11").Value = 0.33
$ws.Range("N11").ClearContents()
$ws.Range("O11").Value = 0.31
$ws.Range("R11").Value = 0.97
$ws.Range("S11").Value = 1.17
$ws.Range("T11").Value = 0.83
$ws.Range("U11").Value = 0.22
$ws.Range("V11").Value = 0.93
